$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 326, shifting existing rows 326:427 down to 327:428.
$ws.Rows("326:326").Insert()

# Populate the new row 326 with the new data record.
$ws.Range("A326").Value = 8
$ws.Range("B326").Value = "Terminal La Palmera de La Serena"
$ws.Range("C326").Value = "Coquimbo"
$ws.Range("D326").Value = 44876
$ws.Range("E326").Value = 4
$ws.Range("F326").Value = 100112032
$ws.Range("G326").Value = "Zapallo italiano"
$ws.Range("H326").Value = "Sin especificar"
$ws.Range("I326").Value = "Primera"
$ws.Range("J326").Value = 520
$ws.Range("K326").Value = 9500
$ws.Range("L326").Value = 10000
$ws.Range("M326").Value = 9750
$ws.Range("N326").Value = "`$/caja 60 unidades"
$ws.Range("O326").Value = "Provincia de Limarí"
$ws.Range("P326").Value = 162
$ws.Range("Q326").Value = 60
$ws.Range("R326").Value = "Hortaliza"
